$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet3 ("MpMatchTriggerType / MpMatchActionType / MpMatcherType" table):
#   - H2:H8 change from a 0..6 sequence to a bit-flag sequence (1,2,4,...,64)
#   - a new row 9 is appended: G9 = "Automatic" (new shared string), H9 = 128
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Sheet3")

$ws3.Range("H2").Value = 1
$ws3.Range("H3").Value = 2
$ws3.Range("H4").Value = 4
$ws3.Range("H5").Value = 8
$ws3.Range("H6").Value = 16
$ws3.Range("H7").Value = 32
$ws3.Range("H8").Value = 64

$ws3.Range("G9").Value = "Automatic"
$ws3.Range("H9").Value = 128

# ---------------------------------------------------------------------------
# Selection / view-state updates (one changed sheetView per sheet)
# ---------------------------------------------------------------------------

# Sheet1: selection moves to E4, view scrolled down so row 4 is at the top
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Activate()
$aw1 = $excel.ActiveWindow
$aw1.ScrollRow = 4
$aw1.ScrollColumn = 1
$ws1.Range("E4").Select()

# Sheet2: selection moves from E4 to D7
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Activate()
$ws2.Range("D7").Select()

# Sheet3: selection moves from J3 to L9 (and stays the active/visible sheet)
$ws3.Activate()
$ws3.Range("L9").Select()
